# Auto-generated: apply scheduled market-price refresh to Golem_Profits sheets
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) for the
# affected leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = -130
$ws.Range("N12").Value = -640

$ws.Range("H33").Value = 247.90475
$ws.Range("I33").Value = 274.44446
$ws.Range("J33").Value = 88.666664
$ws.Range("K33").Value = 274.44446
$ws.Range("L33").Value = 88.666664
$ws.Range("M33").Value = -45.44445999999999
$ws.Range("N33").Value = -546.666664

$ws.Range("H40").Value = 1149.5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1149.5
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1149.5
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1499.5

$ws.Range("H51").Value = 109999.8
$ws.Range("I51").Value = 12500
$ws.Range("K51").Value = 12500
$ws.Range("M51").Value = -12016

$ws.Range("H92").Value = 333333900
$ws.Range("I92").Value = 333333900
$ws.Range("K92").Value = 333333900
$ws.Range("M92").Value = -333332652

$ws.Range("H94").Value = 589
$ws.Range("I94").Value = 589
$ws.Range("K94").Value = 589
$ws.Range("M94").Value = -138

$ws.Range("H132").Value = 1470
$ws.Range("I132").Value = 1403.75
$ws.Range("K132").Value = 4211.25
$ws.Range("M132").Value = -1681.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 355.70587
$ws.Range("I2").Value = 355.70587
$ws.Range("K2").Value = 355.70587
$ws.Range("M2").Value = -242.70587

$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

$ws.Range("H116").Value = 355.70587
$ws.Range("I116").Value = 355.70587
$ws.Range("K116").Value = 355.70587
$ws.Range("M116").Value = 1938.29413

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 355.70587
$ws.Range("I3").Value = 355.70587
$ws.Range("K3").Value = 355.70587
$ws.Range("M3").Value = -241.70587

$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws.Range("H99").Value = 438
$ws.Range("I99").Value = 438
$ws.Range("K99").Value = 438
$ws.Range("M99").Value = 1060

$ws.Range("H107").Value = 25105.21
$ws.Range("I107").Value = 26277.777
$ws.Range("K107").Value = 26277.777
$ws.Range("M107").Value = -24357.777

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1200
$ws.Range("I94").Value = 1200
$ws.Range("K94").Value = 1200
$ws.Range("M94").Value = -749

$ws.Range("H119").Value = 64000
$ws.Range("J119").Value = 64000
$ws.Range("L119").Value = 64000
$ws.Range("N119").Value = -73676

$ws.Range("H132").Value = 1166.6666
$ws.Range("I132").Value = 1166.6666
$ws.Range("K132").Value = 3499.9998
$ws.Range("M132").Value = -969.9998000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1622.2727
$ws.Range("I34").Value = 443.125
$ws.Range("J34").Value = 4766.6665
$ws.Range("K34").Value = 1329.375
$ws.Range("L34").Value = 14299.9995
$ws.Range("M34").Value = -1245.375
$ws.Range("N34").Value = -14467.9995

$ws.Range("H39").Value = 15714.286
$ws.Range("J39").Value = 15714.286
$ws.Range("L39").Value = 47142.858
$ws.Range("N39").Value = -47730.858

$ws.Range("H55").Value = 2081.1875
$ws.Range("J55").Value = 3979.8
$ws.Range("L55").Value = 11939.4
$ws.Range("N55").Value = -12293.4

$ws.Range("H60").Value = 506
$ws.Range("J60").Value = 506
$ws.Range("L60").Value = 1518
$ws.Range("N60").Value = -2020

$ws.Range("H129").Value = 1071.7778
$ws.Range("J129").Value = 1118.1428
$ws.Range("L129").Value = 3354.4284
$ws.Range("N129").Value = -13354.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 8574.75
$ws.Range("J41").Value = 4649.5
$ws.Range("L41").Value = 4649.5
$ws.Range("N41").Value = -5359.5

$ws.Range("H80").Value = 1833.3334
$ws.Range("I80").Value = 1250
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 1250
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -252
$ws.Range("N80").Value = -4996

$ws.Range("H83").Value = 1833.3334
$ws.Range("I83").Value = 1250
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 6250
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -1258
$ws.Range("N83").Value = -24984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2322.4
$ws.Range("I16").Value = 2666.5
$ws.Range("J16").Value = 1806.25
$ws.Range("K16").Value = 2666.5
$ws.Range("L16").Value = 1806.25
$ws.Range("M16").Value = -2496.5
$ws.Range("N16").Value = -2146.25

$ws.Range("H46").Value = 5908.3335
$ws.Range("I46").Value = 4250
$ws.Range("J46").Value = 6737.5
$ws.Range("K46").Value = 4250
$ws.Range("L46").Value = 6737.5
$ws.Range("M46").Value = -4062
$ws.Range("N46").Value = -7113.5

$ws.Range("H68").Value = 1950
$ws.Range("I68").Value = 1960
$ws.Range("J68").Value = 1900
$ws.Range("K68").Value = 1960
$ws.Range("L68").Value = 1900
$ws.Range("M68").Value = -1211
$ws.Range("N68").Value = -3398

$ws.Range("H71").Value = 1950
$ws.Range("I71").Value = 1960
$ws.Range("J71").Value = 1900
$ws.Range("K71").Value = 9800
$ws.Range("L71").Value = 9800
$ws.Range("M71").Value = -6056
$ws.Range("N71").Value = -16988

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
